# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking price strings (column D) are forced to Text format before
# assignment so Excel keeps them as literal strings (matching the source
# data, which stores prices like "2.998.78" / "1.00" as text, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.018.60"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "3.003.50"
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.88"
$ws.Range("E5").Value = "  -4.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.83"
$ws.Range("E6").Value = "  -6.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.002.87"
$ws.Range("E8").Value = "  -5.41%  "
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.135"
$ws.Range("E10").Value = "  -5.81%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").Value = "  -6.00%  "
$ws.Range("E13").Value = "  -5.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.99"
$ws.Range("E14").Value = "  -4.89%  "
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "3.496.05"
$ws.Range("E16").Value = "  -5.42%  "
$ws.Range("D17").Value = "61.085.11"
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").Value = "3.002.55"
$ws.Range("E18").Value = "  -5.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.23"
$ws.Range("E19").Value = "  -5.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "437.78"
$ws.Range("E20").Value = "  -4.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.17"
$ws.Range("E22").Value = "  -6.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.15"
$ws.Range("E23").Value = "  -6.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.10"
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("E25").Value = "  -6.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -6.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  -7.59%  "
$ws.Range("E30").Value = "  -6.87%  "
$ws.Range("E31").Value = "  -9.98%  "
$ws.Range("E32").Value = "  -6.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0941"
$ws.Range("E33").Value = "  -9.43%  "
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.951"
$ws.Range("E35").Value = "  -8.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.55"
$ws.Range("E36").Value = "  -5.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.21"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "0.0₃0672"
$ws.Range("E38").Value = "  -5.81%  "
$ws.Range("E39").Value = "  -6.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.75"
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "374.78"
$ws.Range("E42").Value = "  -6.33%  "
$ws.Range("D43").Value = "2.669.91"
$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("E44").Value = "  -8.83%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -7.11%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.58"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.46"
$ws.Range("E48").Value = "  -5.41%  "
$ws.Range("E49").Value = "  -7.99%  "
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.42"
$ws.Range("E51").Value = "  -8.31%  "
